$wb = $excel.ActiveWorkbook

# Rename the existing sheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "first"

# Add three new sheets after the last existing sheet, in order
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "second"

$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "third"

$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "fourth"

Write-Host "Sheets:"
for ($i=1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Host $wb.Worksheets.Item($i).Name
}
